# "variables get markdown processed"
#
# The title-block placeholders ([TITLE]/by/[AUTHOR]/[PROFESSOR]/etc.) are now
# rendered through markdown processing, which causes Pandoc's OOXML writer to
# start emitting xml:space="preserve" on their <w:t> runs. The professor /
# mnemonic / date line also gets rebuilt so each text segment and each line
# break lives in its own <w:r> run, the date advances a day, and the
# paragraph after the title block switches from the "FirstParagraph" style to
# "BodyText" (keeping its page-break-before setting).
#
# We rebuild the first five paragraphs (title, "by", author, professor block,
# and the paragraph that starts the body) in one shot via Range.InsertXML,
# which lets us control the exact run layout and the xml:space attribute.

$d = $word.ActiveDocument

$firstPara = $d.Paragraphs(1)
$lastPara  = $d.Paragraphs(5)

$firstText = $firstPara.Range.Text.TrimEnd([char]13)
$lastText  = $lastPara.Range.Text.TrimEnd([char]13)
if ($firstText -ne "[TITLE]") {
    throw "edit.ps1: expected paragraph 1 to be '[TITLE]', found '$firstText'"
}
if ($lastText -ne "Nested lists!") {
    throw "edit.ps1: expected paragraph 5 to be 'Nested lists!', found '$lastText'"
}

$rng = $d.Range($firstPara.Range.Start, $lastPara.Range.End)

$xml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="Title"/></w:pPr><w:r><w:t xml:space="preserve">[TITLE]</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Author"/></w:pPr><w:r><w:t xml:space="preserve">by</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Author"/></w:pPr><w:r><w:t xml:space="preserve">[AUTHOR]</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Author"/></w:pPr><w:r><w:t xml:space="preserve">[PROFESSOR]</w:t></w:r><w:r><w:br/></w:r><w:r><w:t xml:space="preserve">[MNEMONIC] &#8212; [CLASS_NAME]</w:t></w:r><w:r><w:br/></w:r><w:r><w:t xml:space="preserve">May 30, 2022</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="BodyText"/><w:pageBreakBefore/></w:pPr><w:r><w:t xml:space="preserve">Nested lists!</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$result = $rng.InsertXML($xml)
